# Refresh the cryptos list (Price / Volume(1h) columns) with the latest
# scraped figures, mirroring the GitHub Actions data-refresh commit.
#
# Note: several "Price" values look like plain numbers (e.g. "310.55").
# Excel auto-converts such literals to numeric cells, which would both
# change the stored cell type and introduce floating-point rounding
# (310.55 -> 310.55000000000001). Prefixing those with a leading
# apostrophe forces Excel to keep them as literal text, exactly matching
# the original formatting (values that already contain 2+ dots, like
# "42.412.44", are never numeric so they don't need the prefix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.412.44"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "2.531.51"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'310.55"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "'98.96"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("D7").Value = "'0.567"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "'35.69"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "'7.33"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "2.921.66"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "'15.73"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D16").Value = "2.516.28"
$ws.Range("E16").Value = "  -3.47%  "
$ws.Range("D17").Value = "'0.824"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "42.418.68"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "'6.79"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "'12.20"
$ws.Range("E21").Value = "  -3.26%  "
$ws.Range("D22").Value = "'69.14"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "'243.81"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'25.97"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("D29").Value = "'39.14"
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "'157.16"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("D32").Value = "'5.71"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").Value = "'2.79"
$ws.Range("E33").Value = "  +15.02%  "
$ws.Range("D34").Value = "'0.0795"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -3.21%  "
$ws.Range("D36").Value = "'2.02"
$ws.Range("E36").Value = "  -5.17%  "
$ws.Range("D37").Value = "'3.16"
$ws.Range("E37").Value = "  -7.56%  "
$ws.Range("D38").Value = "'18.13"
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'4.28"
$ws.Range("E41").Value = "  +9.16%  "
$ws.Range("D42").Value = "'21.72"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'3.30"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "'0.0297"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("D46").Value = "1.962.64"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").Value = "'8.90"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("D48").Value = "2.776.41"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").Value = "'80.90"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").Value = "'0.192"
$ws.Range("D51").Value = "'0.848"
$ws.Range("E51").Value = "  +8.19%  "
